# The "Bibliografia" paragraph originally held all eight references
# concatenated back-to-back with no separation. The edit inserts a blank
# line (two manual line breaks, i.e. <w:br/><w:br/>) between each
# reference, turning it into a readable list while keeping everything in
# a single run/paragraph, exactly as in the target revision.

$d = $word.ActiveDocument

$refs = @(
    "BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gestão Logística da Cadeia de Suprimentos. [s.l.] AMGH, 2013. 472 p.",
    "BARBIERI, J. C. Gestão Ambiental Empresarial: conceitos, modelos e instrumentos. Editora Saraiva, 2004.",
    "ALLEN, D.T.; SHONNARD, D. R., Sustainable Engineering: concepts, design and case studies, Prentice Hall, 2015. ",
    "AKKUCUK, U. Handbook of Research on Sustainable Supply Chain Management for the Global Economy. [s.l.] IGI Global, 2020. 409 p.",
    "BOUCHERY, Y.; CORBETT, C. J.; FRANSOO, J. C.; TAN, T. (ed.). Sustainable Supply Chains. Cham: Springer International Publishing, 2017. v. 4. 130 p.",
    "SCHMIDT, M.; GIOVANNUCCI, D.; PALEKHOV, D.; HANSMANN, B. (ed.). Sustainable Global Value Chains. Cham: Springer International Publishing, 2019. v. 2. 304 p.",
    "LAVE, L. B.; HENDRICKSON, C. T. Environmental Life Cycle Assessment of Goods and Services, Editora John Hopkins, 2006.",
    "LEITE, P. R. Logística Reversa - Meio Ambiente e Competitividade, Editora Prentice Hall: São Paulo, 2002."
)

$oldText = [string]::Join("", $refs)
# "^l" is Word's Find/Replace code for a manual line break (<w:br/>);
# two in a row produce the blank-line separation seen in the diff.
$newText = [string]::Join("^l^l", $refs)

$found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Bibliography paragraph text not found; nothing replaced."
}

Write-Output "found=$found"
